# TC21_Canine_Filter_Breed-Doberman.xlsx
# Fixed variables and query errors in Bread from TC01 to TC30
#
# The Cypher query stored in the "CasesTab" row (cell B2 on sheet "startup")
# referenced a `co` (cohort) variable in its RETURN clause that wasn't meant
# to be part of this query's output. Drop the trailing
# "coalesce(co.cohort_description, '') AS `Cohort`" return column (and its
# now-unneeded trailing blank line) to fix the query.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Doberman Pinscher'] 
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

# Here-strings keep a trailing newline before the closing '@ marker; strip it
# since the new shared-string value has no trailing newline.
$newQuery = $newQuery.TrimEnd("`r", "`n")

$ws.Range("B2").Value = $newQuery

$ws.Range("B2").Select()
